$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (General -> Text) on Price cells whose new value would otherwise
# be auto-parsed as a number by Excel, so they stay literal strings like the original.
$textRows = @(4,5,6,7,8,9,10,11,12,14,15,16,17,18,19,20,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51)
foreach ($r in $textRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Updated Price (D) and Volume(1h) (E) values
$ws.Cells.Item(2, 4).Value = "27.560.29"
$ws.Cells.Item(2, 5).Value = "  -1.28%  "
$ws.Cells.Item(3, 4).Value = "1.845.53"
$ws.Cells.Item(3, 5).Value = "  -2.07%  "
$ws.Cells.Item(4, 4).Value = "1.004"
$ws.Cells.Item(4, 5).Value = "  -1.45%  "
$ws.Cells.Item(5, 4).Value = "333.63"
$ws.Cells.Item(5, 5).Value = "  -0.53%  "
$ws.Cells.Item(6, 4).Value = "1.004"
$ws.Cells.Item(6, 5).Value = "  -1.30%  "
$ws.Cells.Item(7, 4).Value = "0.4645"
$ws.Cells.Item(7, 5).Value = "  -0.83%  "
$ws.Cells.Item(8, 4).Value = "0.3857"
$ws.Cells.Item(8, 5).Value = "  -1.28%  "
$ws.Cells.Item(9, 4).Value = "46.22"
$ws.Cells.Item(9, 5).Value = "  -1.37%  "
$ws.Cells.Item(10, 4).Value = "0.07922"
$ws.Cells.Item(10, 5).Value = "  -0.39%  "
$ws.Cells.Item(11, 4).Value = "0.9966"
$ws.Cells.Item(11, 5).Value = "  -1.58%  "
$ws.Cells.Item(12, 4).Value = "21.53"
$ws.Cells.Item(12, 5).Value = "  -0.88%  "
$ws.Cells.Item(13, 4).Value = "1.846.21"
$ws.Cells.Item(13, 5).Value = "  -2.04%  "
$ws.Cells.Item(14, 4).Value = "5.936"
$ws.Cells.Item(14, 5).Value = "  -0.30%  "
$ws.Cells.Item(15, 4).Value = "7.124"
$ws.Cells.Item(15, 5).Value = "  +0.07%  "
$ws.Cells.Item(16, 4).Value = "1.005"
$ws.Cells.Item(16, 5).Value = "  -1.45%  "
$ws.Cells.Item(17, 4).Value = "89.08"
$ws.Cells.Item(17, 5).Value = "  +1.89%  "
$ws.Cells.Item(18, 4).Value = "0.06647"
$ws.Cells.Item(18, 5).Value = "  -2.03%  "
$ws.Cells.Item(19, 4).Value = "0.00001036"
$ws.Cells.Item(19, 5).Value = "  -1.02%  "
$ws.Cells.Item(20, 4).Value = "17.02"
$ws.Cells.Item(20, 5).Value = "  +0.04%  "
$ws.Cells.Item(21, 5).Value = "  -1.17%  "
$ws.Cells.Item(22, 4).Value = "27.564.64"
$ws.Cells.Item(22, 5).Value = "  -1.26%  "
$ws.Cells.Item(23, 4).Value = "5.388"
$ws.Cells.Item(23, 5).Value = "  -1.47%  "
$ws.Cells.Item(24, 4).Value = "10.91"
$ws.Cells.Item(24, 5).Value = "  -0.18%  "
$ws.Cells.Item(25, 4).Value = "2.305"
$ws.Cells.Item(25, 5).Value = "  -2.37%  "
$ws.Cells.Item(26, 4).Value = "158.60"
$ws.Cells.Item(26, 5).Value = "  -0.77%  "
$ws.Cells.Item(27, 4).Value = "19.53"
$ws.Cells.Item(27, 5).Value = "  -2.35%  "
$ws.Cells.Item(28, 4).Value = "2.107"
$ws.Cells.Item(28, 5).Value = "  +1.23%  "
$ws.Cells.Item(29, 4).Value = "5.415"
$ws.Cells.Item(29, 5).Value = "  -1.04%  "
$ws.Cells.Item(30, 4).Value = "119.98"
$ws.Cells.Item(30, 5).Value = "  -0.91%  "
$ws.Cells.Item(31, 4).Value = "0.9799"
$ws.Cells.Item(31, 5).Value = "  +2.49%  "
$ws.Cells.Item(32, 4).Value = "0.09408"
$ws.Cells.Item(32, 5).Value = "  -1.54%  "
$ws.Cells.Item(33, 4).Value = "3.586"
$ws.Cells.Item(33, 5).Value = "  -1.87%  "
$ws.Cells.Item(34, 4).Value = "5.289"
$ws.Cells.Item(34, 5).Value = "  -0.82%  "
$ws.Cells.Item(35, 4).Value = "1.344"
$ws.Cells.Item(35, 5).Value = "  -0.34%  "
$ws.Cells.Item(36, 4).Value = "0.06038"
$ws.Cells.Item(36, 5).Value = "  -1.18%  "
$ws.Cells.Item(37, 4).Value = "0.02231"
$ws.Cells.Item(37, 5).Value = "  -0.61%  "
$ws.Cells.Item(38, 4).Value = "8.313"
$ws.Cells.Item(38, 5).Value = "  +1.91%  "
$ws.Cells.Item(39, 4).Value = "1.181"
$ws.Cells.Item(39, 5).Value = "  -2.22%  "
$ws.Cells.Item(40, 4).Value = "0.5895"
$ws.Cells.Item(40, 5).Value = "  -0.36%  "
$ws.Cells.Item(41, 4).Value = "0.1868"
$ws.Cells.Item(41, 5).Value = "  -1.30%  "
$ws.Cells.Item(42, 4).Value = "10.32"
$ws.Cells.Item(42, 5).Value = "  +0.93%  "
$ws.Cells.Item(43, 4).Value = "1.241"
$ws.Cells.Item(43, 5).Value = "  -2.27%  "
$ws.Cells.Item(44, 4).Value = "0.5591"
$ws.Cells.Item(44, 5).Value = "  -1.13%  "
$ws.Cells.Item(45, 4).Value = "12.14"
$ws.Cells.Item(45, 5).Value = "  -0.39%  "
$ws.Cells.Item(46, 4).Value = "1.909"
$ws.Cells.Item(46, 5).Value = "  -0.96%  "
$ws.Cells.Item(47, 4).Value = "0.06694"
$ws.Cells.Item(47, 5).Value = "  -2.37%  "
$ws.Cells.Item(48, 4).Value = "111.09"
$ws.Cells.Item(48, 5).Value = "  -2.48%  "
$ws.Cells.Item(49, 4).Value = "1.053"
$ws.Cells.Item(49, 5).Value = "  -1.12%  "
$ws.Cells.Item(50, 4).Value = "1.003"
$ws.Cells.Item(50, 5).Value = "  -1.48%  "
$ws.Cells.Item(51, 4).Value = "70.17"
$ws.Cells.Item(51, 5).Value = "  -0.77%  "
